$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-point existing per-cell styles before any values move, so we can
#    reuse already-defined cellXfs (style indices) rather than create new
#    ones.  This mirrors the final layout where:
#      - C2/C3 and G2/G3 use the "bold-ish black font + full border" style
#        (previously used only by E2/E3)
#      - D2/D3 use a brand new style: same font, but border without a left
#        edge (so it doesn't double up against C's right edge)
#      - E2/E3 fall back to the plain bordered style (like B2/B3)
# ---------------------------------------------------------------------------

# C2/C3 <- take current E2/E3 formatting (font+border) before it's reset
$ws.Range("E2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null

# E2:F2 / E3:F3 <- reset back to the plain bordered look (copy from B2/B3)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2:F2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("E3:F3").PasteSpecial(-4122) | Out-Null

# G2/G3 <- same styling as C2/C3
$ws.Range("C2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null

# D2:D3 <- new style: black font + border with no left edge
$ws.Range("D2:D3").Font.Color = 0
$ws.Range("D2:D3").Borders.Item(7).LineStyle = -4142

# N1 header <- same look as the other header cells
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null

# N2/N3 <- plain bordered look like the rest of the data cells
$ws.Range("B2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Update header row text (columns E-J were re-ordered, N is new)
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "Tier"
$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"
$ws.Range("G1").Value = "Client"
$ws.Range("H1").Value = "Lob"
$ws.Range("J1").Value = "Product Name"

# ---------------------------------------------------------------------------
# 3) Update data rows 2 and 3
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Qu18-001"
$ws.Range("C2").Value = "SIPL6118"
$ws.Range("D2").Value = "SIPL4167"
$ws.Range("E2").Value = "SIPL5317"
$ws.Range("F2").Value = "SIPL5317"
$ws.Range("G2").Value = "Qualia"
$ws.Range("H2").Value = "Title"
$ws.Range("J2").Value = "Current Owner Search"
$ws.Range("K2").Value = "AL"
$ws.Range("L2").Value = "Shelby"

$ws.Range("B3").Value = "Qu18-002"
$ws.Range("C3").Value = "SIPL5316"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"
$ws.Range("G3").Value = "Qualia"
$ws.Range("H3").Value = "Title"
$ws.Range("J3").Value = "Full Search"

# ---------------------------------------------------------------------------
# 4) Column widths - best effort match of the new layout
# ---------------------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 20.33203125
$ws.Columns("E:E").ColumnWidth = 8.43
$ws.Columns("G:G").ColumnWidth = 6.21875
$ws.Columns("H:H").ColumnWidth = 4.33203125
$ws.Columns("I:I").ColumnWidth = 13.6640625
$ws.Columns("J:J").ColumnWidth = 19
$ws.Columns("N:N").ColumnWidth = 11.5546875

# ---------------------------------------------------------------------------
# 5) Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("E5").Select() | Out-Null
